# Add five new vocabulary entries (severity, sensation, smuggle, slope, soak)
# to the end of the word list on Sheet1, matching the layout/styling of the
# existing rows (word | synonym/definition | example1 | example2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the direct formatting (font/alignment/etc, style index "2") from the
# last existing data row (B102:D102, which already uses that formatting) onto
# the five new rows, covering columns A:D.
$ws.Range("B102:D102").Copy() | Out-Null
$ws.Range("A103:D107").PasteSpecial(-4122) | Out-Null

# Row 103 - severity
$ws.Range("A103").Value = "severity"
$ws.Range("B103").Value = "severe problems, injuries, illnesses etc are very bad or very serious"
$ws.Range("C103").Value = "he risk and severity of sunborn depend on he body's natural skin color."
$ws.Range("D103").Value = "His injuries were quite severe."
$ws.Rows.Item(103).RowHeight = 75

# Row 104 - sensation
$ws.Range("A104").Value = "sensation"
$ws.Range("B104").Value = "a feeling that you get from one of your five senses, especially the sense of touch"
$ws.Range("C104").Value = "I experienced no sensation in my left foot."
$ws.Range("D104").Value = "One sign of a heart attack is a tingling sensation in the left arm."
$ws.Rows.Item(104).RowHeight = 60

# Row 105 - smuggle
$ws.Range("A105").Value = "smuggle"
$ws.Range("B105").Value = "to take something or someone illegally from one country to another"
$ws.Range("C105").Value = "if you try to smuggle drug you are stupid."
$ws.Range("D105").Value = "The guns were smuggled across the border."
$ws.Rows.Item(105).RowHeight = 45

# Row 106 - slope
$ws.Range("A106").Value = "slope"
$ws.Range("B106").Value = "a surface of which one end or side is at a higher level than another; a rising or falling surface."
$ws.Range("C106").Value = "the house builders slopped the roof..."
$ws.Range("D106").Value = "the roof should have a slope sufficient for proper drainage"
$ws.Rows.Item(106).RowHeight = 75

# Row 107 - soak
$ws.Range("A107").Value = "soak"
$ws.Range("B107").Value = "if you soak something, or if you let it soak, you keep it covered with a liquid for a period of time, especially in order to make it softer or easier to clean"
$ws.Range("C107").Value = "Soak the clothes in cold water."
$ws.Range("D107").Value = "soak the beans overnight in water"
$ws.Rows.Item(107).RowHeight = 105

# Move the visible selection to match where the author ended up (E107).
$ws.Range("E107").Select() | Out-Null

Write-Output "done"
